# "etiquetando e alterando notebook"
# Re-labels several tweets in column B (the "Etiquetas" column) of Sheet1,
# then leaves the selection the way the author left it (rows 2:301 selected,
# scrolled down near the bottom of the sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Re-labelled tweets (column B, "Etiquetas") ---
$updates = @{
    22  = 5
    24  = 5
    25  = 5
    26  = 5
    27  = 5
    29  = 5
    35  = 5
    36  = 5
    37  = 5
    45  = 5
    47  = 5
    48  = 4
    51  = 5
    55  = 1
    56  = 5
    57  = 5
    59  = 5
    61  = 5
    63  = 5
    64  = 5
    66  = 5
    67  = 5
    68  = 5
    69  = 5
    71  = 5
    74  = 5
    75  = 5
    77  = 5
    79  = 5
    81  = 5
    82  = 5
    93  = 4
    98  = 5
    99  = 5
    101 = 5
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 2).Value = $updates[$row]
}

# --- Final selection / scroll state, matching the author's saved view ---
$ws.Range("A2:XFD301").Select()
try { $excel.ActiveWindow.ScrollRow = 279 } catch {}
try { $excel.ActiveWindow.ScrollColumn = 1 } catch {}
